$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update F2:F9 time_taken values on the "data" sheet
$ws1.Range("F2").Value = "2021-10-05 14:21:33.103478"
$ws1.Range("F3").Value = "2021-10-05 14:21:33.103486"
$ws1.Range("F4").Value = "2021-10-05 14:21:33.103489"
$ws1.Range("F5").Value = "2021-10-05 14:21:33.103492"
$ws1.Range("F6").Value = "2021-10-05 14:21:33.103495"
$ws1.Range("F7").Value = "2021-10-05 14:21:33.103497"
$ws1.Range("F8").Value = "2021-10-05 14:21:33.103500"
$ws1.Range("F9").Value = "2021-10-05 14:21:33.103503"

# Add a new "metadata" worksheet after "data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Set header row values first (order matters: format copy must come after value set)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Set data row values
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Mitochondrial disorder with complex II deficiency"
$ws2.Range("C2").Value = 535
$ws2.Range("E2").Value = "2020-02-17T15:54:32.903706Z"
$ws2.Range("F2").Value = "2021-10-05 14:21:33.099719"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/535/?format=json"

# "1.3" must stay a text value (not get auto-coerced to a number). Use a
# scratch cell holding a text-producing formula, then paste its value only.
$scratch = $ws2.Range("Z100")
$scratch.Formula = '="1.3"'
$scratch.Copy()
$ws2.Range("D2").PasteSpecial(-4104)
$scratch.Clear()

# Copy header formatting (bold, bordered style) from "data" sheet header row onto metadata header row
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

# Copy the A-column style (used for numeric index column) from "data" sheet onto metadata's A2
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
